$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

Set-TextValue "D2" "327.86"
Set-TextValue "E2" "1.37%"
Set-TextValue "B3" "OKB"
Set-TextValue "C3" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D3" "43.99"
Set-TextValue "E3" "-1.06%"
Set-TextValue "B4" "HuobiToken"
Set-TextValue "C4" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D4" "5.500"
Set-TextValue "E4" "-0.06%"
Set-TextValue "B5" "Cronos"
Set-TextValue "C5" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D5" "0.08013"
Set-TextValue "E5" "-0.27%"
Set-TextValue "B6" "FTXToken"
Set-TextValue "C6" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D6" "2.017"
Set-TextValue "E6" "7.06%"
Set-TextValue "B7" "GateToken"
Set-TextValue "C7" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "4.315"
Set-TextValue "E7" "-0.51%"
Set-TextValue "B8" "BTSEToken"
Set-TextValue "C8" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D8" "2.584"
Set-TextValue "E8" "-4.00%"
Set-TextValue "B9" "MXToken"
Set-TextValue "C9" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D9" "0.9508"
Set-TextValue "E9" "1.00%"
Set-TextValue "B10" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D10" "0.1122"
Set-TextValue "E10" "-4.62%"
Set-TextValue "B11" "WazirX"
Set-TextValue "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1872"
Set-TextValue "E11" "-0.23%"
Set-TextValue "B12" "MCDex"
Set-TextValue "C12" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D12" "10.63"
Set-TextValue "E12" "25.52%"
Set-TextValue "B13" "MandalaExchangeToken"
Set-TextValue "C13" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D13" "0.09887"
Set-TextValue "E13" "-0.40%"
Set-TextValue "B14" "BitrueCoin"
Set-TextValue "C14" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.04582"
Set-TextValue "E14" "10.01%"
Set-TextValue "B15" "BitMartToken"
Set-TextValue "C15" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.1067"
Set-TextValue "E15" "0.30%"
Set-TextValue "B16" "BitForexToken"
Set-TextValue "C16" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D16" "0.001264"
Set-TextValue "E16" "-0.49%"
Set-TextValue "B17" "CoinExToken"
Set-TextValue "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D17" "0.04069"
Set-TextValue "E17" "-4.32%"
Set-TextValue "B18" "TigerCash"
Set-TextValue "C18" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D18" "0.005932"
Set-TextValue "E18" "0.82%"
Set-TextValue "E19" "-6.65%"
Set-TextValue "E20" "-0.30%"
Set-TextValue "D21" "0.1406"
Set-TextValue "E21" "2.26%"
Set-TextValue "E22" "0.17%"
Set-TextValue "D23" "0.001259"
Set-TextValue "E23" "1.32%"
Set-TextValue "D24" "0.004332"
Set-TextValue "E25" "-6.36%"
Set-TextValue "D26" "0.0003740"
Set-TextValue "E26" "-6.69%"
Set-TextValue "D38" "0.02563"
Set-TextValue "E38" "-2.52%"
Set-TextValue "D39" "0.05679"
Set-TextValue "E39" "3.53%"
Set-TextValue "D40" "0.007530"
Set-TextValue "E40" "-2.59%"
Set-TextValue "D41" "0.1396"
Set-TextValue "E41" "0.50%"
Set-TextValue "D42" "0.007602"
Set-TextValue "E42" "11.94%"
Set-TextValue "D43" "0.002018"
Set-TextValue "E43" "-1.98%"
Set-TextValue "D44" "0.008875"
Set-TextValue "E44" "-3.85%"
Set-TextValue "D45" "0.00007104"
Set-TextValue "E45" "-0.45%"
Set-TextValue "E46" "-0.69%"
Set-TextValue "E47" "54.73%"
Set-TextValue "D48" "0.003099"
Set-TextValue "E48" "-8.94%"
Set-TextValue "D49" "0.00002098"
Set-TextValue "E49" "-0.69%"
Set-TextValue "D50" "0.0001998"
Set-TextValue "E50" "-0.69%"
